$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("M2").Value = 1.442875
$ws.Range("N2").Value = 4.328625
$ws.Range("O2").Value = 0.02047893724893121
$ws.Range("P2").Value = 0.02047893724893121
$ws.Range("Q2").Value = 30.56874926904167
$ws.Range("R2").Value = 275.118743421375
$ws.Range("S2").Value = 0.001795007848108377
$ws.Range("T2").Value = 0.001795007848108377
$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("O3").Value = 0.1473796107804731
$ws.Range("P3").Value = 0.1473796107804731
$ws.Range("Q3").Value = 219.9923909407144
$ws.Range("R3").Value = 1979.931518466429
$ws.Range("S3").Value = 0.01291803157490088
$ws.Range("T3").Value = 0.01291803157490088
$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 27.934719
$ws.Range("N4").Value = 83.804157
$ws.Range("O4").Value = 0.3964815784233052
$ws.Range("P4").Value = 0.3964815784233051
$ws.Range("Q4").Value = 591.824947422427
$ws.Range("R4").Value = 5326.424526801843
$ws.Range("S4").Value = 0.03475217176796479
$ws.Range("T4").Value = 0.03475217176796479
$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("M5").Value = 30.695086
$ws.Range("N5").Value = 92.085258
$ws.Range("O5").Value = 0.4356598735472906
$ws.Range("P5").Value = 0.4356598735472905
$ws.Range("Q5").Value = 650.3060817643047
$ws.Range("R5").Value = 5852.754735878742
$ws.Range("S5").Value = 0.03818620481217123
$ws.Range("T5").Value = 0.03818620481217123
$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("M6").Value = 1.442875
$ws.Range("N6").Value = 4.328625
$ws.Range("O6").Value = 0.02047893724893121
$ws.Range("P6").Value = 0.02047893724893121
$ws.Range("Q6").Value = 175.7723200254584
$ws.Range("R6").Value = 1581.950880229125
$ws.Range("S6").Value = 0.01032141325603559
$ws.Range("T6").Value = 0.01032141325603559
$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("O7").Value = 0.1473796107804731
$ws.Range("P7").Value = 0.1473796107804731
$ws.Range("S7").Value = 0.07427953169094897
$ws.Range("T7").Value = 0.07427953169094896
$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 27.934719
$ws.Range("N8").Value = 83.804157
$ws.Range("O8").Value = 0.3964815784233052
$ws.Range("P8").Value = 0.3964815784233051
$ws.Range("Q8").Value = 3403.032395660921
$ws.Range("R8").Value = 30627.2915609483
$ws.Range("S8").Value = 0.199827274705175
$ws.Range("T8").Value = 0.199827274705175
$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("M9").Value = 30.695086
$ws.Range("N9").Value = 92.085258
$ws.Range("O9").Value = 0.4356598735472906
$ws.Range("P9").Value = 0.4356598735472905
$ws.Range("Q9").Value = 3739.302766768408
$ws.Range("R9").Value = 33653.72490091567
$ws.Range("S9").Value = 0.2195731907029733
$ws.Range("T9").Value = 0.2195731907029732
$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("M10").Value = 1.442875
$ws.Range("N10").Value = 4.328625
$ws.Range("O10").Value = 0.02047893724893121
$ws.Range("P10").Value = 0.02047893724893121
$ws.Range("Q10").Value = 53.68531320920832
$ws.Range("R10").Value = 483.1678188828749
$ws.Range("S10").Value = 0.003152420718641536
$ws.Range("T10").Value = 0.003152420718641536
$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("O11").Value = 0.1473796107804731
$ws.Range("P11").Value = 0.1473796107804731
$ws.Range("Q11").Value = 386.3540607222596
$ws.Range("R11").Value = 3477.186546500337
$ws.Range("S11").Value = 0.02268684809578858
$ws.Range("T11").Value = 0.02268684809578858
$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 27.934719
$ws.Range("N12").Value = 83.804157
$ws.Range("O12").Value = 0.3964815784233052
$ws.Range("P12").Value = 0.3964815784233051
$ws.Range("Q12").Value = 1039.372183263431
$ws.Range("R12").Value = 9354.349649370879
$ws.Range("S12").Value = 0.06103230490862298
$ws.Range("T12").Value = 0.06103230490862298
$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("M13").Value = 30.695086
$ws.Range("N13").Value = 92.085258
$ws.Range("O13").Value = 0.4356598735472906
$ws.Range("P13").Value = 0.4356598735472905
$ws.Range("Q13").Value = 1142.077661539347
$ws.Range("R13").Value = 10278.69895385413
$ws.Range("S13").Value = 0.06706320718488001
$ws.Range("T13").Value = 0.06706320718488003
$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("M14").Value = 1.442875
$ws.Range("N14").Value = 4.328625
$ws.Range("O14").Value = 0.02047893724893121
$ws.Range("P14").Value = 0.02047893724893121
$ws.Range("Q14").Value = 88.72724479587498
$ws.Range("R14").Value = 798.5452031628748
$ws.Range("S14").Value = 0.005210095426145702
$ws.Range("T14").Value = 0.005210095426145702
$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("O15").Value = 0.1473796107804731
$ws.Range("P15").Value = 0.1473796107804731
$ws.Range("Q15").Value = 638.5383501442329
$ws.Range("R15").Value = 5746.845151298096
$ws.Range("S15").Value = 0.03749519941883465
$ws.Range("T15").Value = 0.03749519941883464
$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 27.934719
$ws.Range("N16").Value = 83.804157
$ws.Range("O16").Value = 0.3964815784233052
$ws.Range("P16").Value = 0.3964815784233051
$ws.Range("Q16").Value = 1717.799983378311
$ws.Range("R16").Value = 15460.1998504048
$ws.Range("S16").Value = 0.1008698270415424
$ws.Range("T16").Value = 0.1008698270415424
$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("M17").Value = 30.695086
$ws.Range("N17").Value = 92.085258
$ws.Range("O17").Value = 0.4356598735472906
$ws.Range("P17").Value = 0.4356598735472905
$ws.Range("Q17").Value = 1887.544249884734
$ws.Range("R17").Value = 16987.8982489626
$ws.Range("S17").Value = 0.110837270847266
$ws.Range("T17").Value = 0.110837270847266
